# Updated cryptos list on Mon Nov  4 15:32:43 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D hold price text that LOOKS like a plain number (no
# thousands separators made of extra dots). Excel's normal typed-value
# coercion would silently turn these into numeric cells (and would also
# drop formatting such as trailing zeros), so force them to remain text
# by marking the cell format as Text ("@") before writing the value -
# matching how these cells were already stored (inline/shared strings).
$textCells = @(
  "D5","D6","D8","D9","D14","D15","D16","D17","D18","D19",
  "D22","D24","D26","D28","D31","D32","D34",
  "D36","D37","D38","D39","D41","D42","D43","D45","D46","D49","D50"
)
foreach ($addr in $textCells) {
  $ws.Range($addr).NumberFormat = "@"
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "68.449.94"
$ws.Range("E2").Value = "  +1.23%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "2.446.42"
$ws.Range("E3").Value = "  +0.96%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.04%  "

# Row 5 - BNB
$ws.Range("D5").Value = "557.83"
$ws.Range("E5").Value = "  +1.29%  "

# Row 6 - Solana
$ws.Range("D6").Value = "162.56"
$ws.Range("E6").Value = "  +2.59%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.01%  "

# Row 8 - XRP
$ws.Range("D8").Value = "0.511"
$ws.Range("E8").Value = "  +3.23%  "

# Row 9 - Dogecoin
$ws.Range("D9").Value = "0.157"
$ws.Range("E9").Value = "  +8.96%  "

# Row 10 - TRON
$ws.Range("E10").Value = "  +0.47%  "

# Row 11 - Toncoin
$ws.Range("E11").Value = "  +2.49%  "

# Row 12 - Cardano
$ws.Range("E12").Value = "  -1.52%  "

# Row 13 - WrappedBTC
$ws.Range("D13").Value = "68.348.24"
$ws.Range("E13").Value = "  +0.93%  "

# Row 14 - ShibaInu
$ws.Range("D14").Value = "0.0000170"
$ws.Range("E14").Value = "  +3.68%  "

# Row 15 - Avalanche
$ws.Range("D15").Value = "23.32"
$ws.Range("E15").Value = "  +2.24%  "

# Row 16 - Chainlink
$ws.Range("D16").Value = "10.44"
$ws.Range("E16").Value = "  -1.65%  "

# Row 17 - BitcoinCash
$ws.Range("D17").Value = "337.06"
$ws.Range("E17").Value = "  +0.22%  "

# Row 18 - Uniswap
$ws.Range("D18").Value = "6.90"
$ws.Range("E18").Value = "  -0.49%  "

# Row 19 - Polkadot
$ws.Range("D19").Value = "3.79"
$ws.Range("E19").Value = "  +2.36%  "

# Row 20 - SuiNetwork
$ws.Range("E20").Value = "  +4.61%  "

# Row 21 - Dai
$ws.Range("E21").Value = "  +0.27%  "

# Row 22 - Litecoin
$ws.Range("D22").Value = "66.77"
$ws.Range("E22").Value = "  +1.58%  "

# Row 23 - NEARProtocol
$ws.Range("E23").Value = "  +2.46%  "

# Row 24 - Aptos
$ws.Range("D24").Value = "8.18"
$ws.Range("E24").Value = "  +3.00%  "

# Row 25 - PEPE
$ws.Range("D25").Value = "0.0₃0817"
$ws.Range("E25").Value = "  +2.22%  "

# Row 26 - InternetComputer(DFINITY)
$ws.Range("D26").Value = "7.20"
$ws.Range("E26").Value = "  +3.43%  "

# Row 27 - FirstDigitalUSD
$ws.Range("E27").Value = "  +0.08%  "

# Row 28 - Bittensor
$ws.Range("D28").Value = "426.01"
$ws.Range("E28").Value = "  +1.91%  "

# Row 29 - Fetch.AI
$ws.Range("E29").Value = "  +3.46%  "

# Row 30 - PancakeSwap
$ws.Range("E30").Value = "  +0.83%  "

# Row 31 - Monero
$ws.Range("D31").Value = "160.52"
$ws.Range("E31").Value = "  +1.93%  "

# Row 32 - WhiteBITCoin
$ws.Range("D32").Value = "19.00"
$ws.Range("E32").Value = "  +0.18%  "

# Row 33
$ws.Range("E33").Value = "  +0.18%  "

# Row 34 - EthereumClassic
$ws.Range("D34").Value = "17.80"
$ws.Range("E34").Value = "  +1.60%  "

# Row 35 - Kaspa
$ws.Range("E35").Value = "  -1.48%  "

# Rows 36 and 37 swap content: PolygonEcosystemToken <-> RenderToken
$ws.Range("B36").Value = "RenderToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D36").Value = "4.38"
$ws.Range("E36").Value = "  +3.05%  "

$ws.Range("B37").Value = "PolygonEcosystemToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D37").Value = "0.296"
$ws.Range("E37").Value = "  -0.01%  "

# Row 38 - Stacks
$ws.Range("D38").Value = "1.47"
$ws.Range("E38").Value = "  +3.55%  "

# Row 39 - ImmutableX
$ws.Range("D39").Value = "1.07"
$ws.Range("E39").Value = "  +1.13%  "

# Row 40 - dogwifhat
$ws.Range("E40").Value = "  +2.48%  "

# Row 41 - Filecoin
$ws.Range("D41").Value = "3.36"
$ws.Range("E41").Value = "  +2.90%  "

# Row 42 - Aave
$ws.Range("D42").Value = "129.78"
$ws.Range("E42").Value = "  -1.55%  "

# Row 43 - Cronos
$ws.Range("D43").Value = "0.0717"
$ws.Range("E43").Value = "  +1.52%  "

# Row 44 - ARBITRUM
$ws.Range("E44").Value = "  +2.97%  "

# Row 45 - Mantle
$ws.Range("D45").Value = "0.562"
$ws.Range("E45").Value = "  +2.34%  "

# Row 46 - Stellar
$ws.Range("D46").Value = "0.0918"
$ws.Range("E46").Value = "  +2.20%  "

# Row 47 - BitgetToken
$ws.Range("E47").Value = "  +1.23%  "

# Row 48 - Optimism
$ws.Range("E48").Value = "  -0.72%  "

# Rows 49 and 50 swap content: InjectiveProtocol <-> THORChain
$ws.Range("B49").Value = "THORChain"
$ws.Range("C49").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D49").Value = "4.92"
$ws.Range("E49").Value = "  -2.59%  "

$ws.Range("B50").Value = "InjectiveProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D50").Value = "16.74"
$ws.Range("E50").Value = "  +2.42%  "

# Row 51 - BabyDogeCoin
$ws.Range("E51").Value = "  +5.65%  "
